# Append the new profit row for 2025-09-23 to the bottom of the data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCell = $ws.Cells.Item(37, 1)
# Force the date to be stored as literal text (matching the existing rows),
# rather than letting Excel auto-convert the "MM/DD/YYYY" string into a date
# serial number.
$dateCell.NumberFormat = "@"
$dateCell.Value = "09/23/2025"
# Drop back to the default/normal style so no explicit style index lingers
# on the cell (matches the formatting of the other data rows).
$dateCell.Style = "Normal"

$ws.Cells.Item(37, 2).Value = 15360.44
